$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = -0.0122
$ws.Range("E2").Value = 0.0707
$ws.Range("K2").Value = -4.199999999999999
$ws.Range("L2").Value = -0.01871657754010695
$ws.Range("M2").Value = 17.4726
$ws.Range("N2").Value = 0.04373617021276596
$ws.Range("O2").Value = -4.160142857142858
$ws.Range("P2").Value = 17.4726
$ws.Range("Q2").Value = 0.04373617021276596
$ws.Range("R2").Value = -4.160142857142858
$ws.Range("U2").Value = 63.9
$ws.Range("V2").Value = 0.1599499374217772
$ws.Range("W2").Value = 0.05583195135433941
$ws.Range("X2").Value = 0.1107539613619843
$ws.Range("Y2").Value = -0.05492201000764491
$ws.Range("Z2").Value = 0.1394395078605605
$ws.Range("AB2").Value = 0.0493736782541972
$ws.Range("AC2").Value = -0.0493736782541972
$ws.Range("AD2").Value = 1662.5
$ws.Range("AF2").Value = 1662.5
$ws.Range("AG2").Value = 1598.6
$ws.Range("AH2").Value = 0.8062560620756547
$ws.Range("AI2").Value = 0.8015524805939926
$ws.Range("AJ2").Value = 0.8000600570542015
$ws.Range("AK2").Value = 0.7952442543030545
$ws.Range("K3").Value = 14.9
$ws.Range("L3").Value = 0.1712643678160919
$ws.Range("M3").Value = 10.2357
$ws.Range("N3").Value = 0.07493191800878478
$ws.Range("O3").Value = 0.6869597315436241
$ws.Range("P3").Value = 10.2357
$ws.Range("Q3").Value = 0.07493191800878478
$ws.Range("R3").Value = 0.6869597315436241
$ws.Range("W3").Value = 0.1808252427184466
$ws.Range("X3").Value = 0.07187194681592372
$ws.Range("Y3").Value = 0.1089532959025229
$ws.Range("Z3").Value = 0.3361669242658423
$ws.Range("AB3").Value = 0.04594278262315273
$ws.Range("AC3").Value = -0.04594278262315273
$ws.Range("AD3").Value = 213.8
$ws.Range("AF3").Value = 213.8
$ws.Range("AG3").Value = 213.8
$ws.Range("AH3").Value = 0.6101598173515983
$ws.Range("AI3").Value = 0.6910148674854558
$ws.Range("AJ3").Value = 0.6101598173515983
$ws.Range("AK3").Value = 0.6910148674854558
$ws.Range("D4").Value = 0.0414
$ws.Range("E4").Value = 0.0707
$ws.Range("K4").Value = 10.1
$ws.Range("L4").Value = 0.1276864728192162
$ws.Range("M4").Value = 7.236899999999999
$ws.Range("N4").Value = 0.04242028135990621
$ws.Range("O4").Value = 0.7165247524752475
$ws.Range("P4").Value = 7.236899999999999
$ws.Range("Q4").Value = 0.04242028135990621
$ws.Range("R4").Value = 0.7165247524752475
$ws.Range("U4").Value = 62.5
$ws.Range("V4").Value = 0.3663540445486518
$ws.Range("W4").Value = 0.05583195135433941
$ws.Range("X4").Value = 0.1828852038932544
$ws.Range("Y4").Value = -0.127053252538915
$ws.Range("Z4").Value = 0.0693555458132398
$ws.Range("AB4").Value = 0.0493736782541972
$ws.Range("AC4").Value = -0.0493736782541972
$ws.Range("AD4").Value = 1139
$ws.Range("AF4").Value = 1139
$ws.Range("AG4").Value = 1076.5
$ws.Range("AH4").Value = 0.8697312156383629
$ws.Range("AI4").Value = 0.8363315955650193
$ws.Range("AJ4").Value = 0.8632026301018363
$ws.Range("AK4").Value = 0.8284592889025704
$ws.Range("D5").Value = -0.0658
$ws.Range("K5").Value = -29.2
$ws.Range("L5").Value = -0.5008576329331046
$ws.Range("U5").Value = 1.4
$ws.Range("V5").Value = 0.01516793066088841
$ws.Range("W5").Value = -0.4634920634920635
$ws.Range("X5").Value = 0.1107539613619843
$ws.Range("Y5").Value = -0.5742460248540477
$ws.Range("Z5").Value = 0.2776190476190476
$ws.Range("AB5").Value = 0.0519250539411089
$ws.Range("AC5").Value = -0.0519250539411089
$ws.Range("AD5").Value = 309.7
$ws.Range("AF5").Value = 309.7
$ws.Range("AG5").Value = 308.3
$ws.Range("AH5").Value = 0.7703980099502488
$ws.Range("AI5").Value = 0.768867924528302
$ws.Range("AJ5").Value = 0.7695956065901148
$ws.Range("AK5").Value = 0.7680617837568511
